$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "height"
$ws.Range("G1").Value = "weight"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.5    # E column -> height value
    $ws.Cells.Item($r, 6).Value = 257    # F column -> weight value
    $ws.Cells.Item($r, 7).Value = 0      # G column -> old fantasy points value
}
